# Update automatico via Actualizar 12-01-2020 13-02-05
#
# Appends the newly reported daily COVID-19 figures for Guatemala
# (8/11/2020 - 29/11/2020) as rows 226-247 of the "Condicion_Pacientes"
# table on sheet "Hoja1", then resizes the table/AutoFilter to cover the
# new range and updates the view position/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Hoja1")
$lo = $ws.ListObjects.Item(1)

$data = @(
    @(8,  11, 2020, 648,  98),
    @(9,  11, 2020, 4894, 769),
    @(10, 11, 2020, 4553, 682),
    @(11, 11, 2020, 5377, 732),
    @(12, 11, 2020, 4532, 580),
    @(13, 11, 2020, 5102, 596),
    @(14, 11, 2020, 1278, 166),
    @(15, 11, 2020, 1096, 147),
    @(16, 11, 2020, 5210, 698),
    @(17, 11, 2020, 4325, 651),
    @(18, 11, 2020, 5158, 685),
    @(19, 11, 2020, 4945, 691),
    @(20, 11, 2020, 4419, 660),
    @(21, 11, 2020, 1338, 212),
    @(22, 11, 2020, 665,  93),
    @(23, 11, 2020, 519,  627),
    @(24, 11, 2020, 4606, 640),
    @(25, 11, 2020, 4617, 696),
    @(26, 11, 2020, 3921, 447),
    @(27, 11, 2020, 4115, 666),
    @(28, 11, 2020, 1580, 173),
    @(29, 11, 2020, 625,  91)
)

$firstNewRow = 226
$lastOldRow  = 225
$lastNewRow  = $firstNewRow + $data.Count - 1   # 247

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $firstNewRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Formula = "=+Condicion_Pacientes[[#This Row],[día]]&""/""&Condicion_Pacientes[[#This Row],[mes]]&""/""&Condicion_Pacientes[[#This Row],[año]]"
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]

    # Match the formatting already used for the preceding data rows
    # (col A: centered + grey fill, cols B-D: centered, E-F: default)
    # by copying the format from the last pre-existing row.
    $ws.Range("A$lastOldRow`:F$lastOldRow").Copy()
    $ws.Range("A$row`:F$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Grow the table (and its AutoFilter) to cover the appended rows.
$lo.Resize($ws.Range("A1:I$lastNewRow"))

# Reflect the view state captured in the edited workbook.
$ws.Application.ActiveWindow.ScrollRow = 222
$ws.Range("F249").Select()
